$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps plain-text values such as '26.806.94' or
# '1.001' instead of Excel reinterpreting them as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.806.94"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "1.870.96"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "301.26"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.5340"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("D8").Value = "0.3744"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "0.07190"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "0.8892"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "0.08165"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "1.885.36"
$ws.Range("E13").Value = "  +26.15%  "
$ws.Range("D14").Value = "92.88"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "5.308"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "14.85"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "0.000008504"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("D20").Value = "26.835.34"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "4.989"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "6.383"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "2.313"
$ws.Range("D25").Value = "146.05"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").Value = "1.731"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "18.03"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "114.01"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "4.723"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "4.632"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").Value = "0.09153"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "0.8043"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").Value = "0.05029"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "1.174"
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("D35").Value = "2.945"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "0.6145"
$ws.Range("E36").Value = "  +6.55%  "
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").Value = "3.192"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").Value = "1.065"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("B41").Value = "Decentraland"
$ws.Range("C41").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D41").Value = "0.5247"
$ws.Range("E41").Value = "  +6.79%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.531"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "8.792"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").Value = "114.81"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").Value = "0.1492"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.650"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").Value = "9.944"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "37.63"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").Value = "0.06054"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "62.28"
$ws.Range("E51").Value = "  -3.30%  "
